$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows appended to the Dictionary sheet describing the "Caltrain Pilot"
# onoff / transfer-sequence variables (rows 546-561).
# Columns: A=Operator, B=Survey_Variable, C=Survey_Response, D=Generic_Variable, E=Generic_Response
# Cells are written row-by-row, left-to-right, in the same order they were
# authored so new shared-string entries are interned in the matching order.

$ws.Cells.Item(546, 1).Value = "Caltrain Pilot"
$ws.Cells.Item(546, 2).Value = "ENTER_STATION"
$ws.Cells.Item(546, 3).Value = "NONCATEGORICAL"
$ws.Cells.Item(546, 4).Value = "onoff_enter_station"
$ws.Cells.Item(546, 5).Value = "NONCATEGORICAL"

$ws.Cells.Item(547, 1).Value = "Caltrain Pilot"
$ws.Cells.Item(547, 2).Value = "EXIT_STATION"
$ws.Cells.Item(547, 3).Value = "NONCATEGORICAL"
$ws.Cells.Item(547, 4).Value = "onoff_exit_station"
$ws.Cells.Item(547, 5).Value = "NONCATEGORICAL"

$ws.Cells.Item(548, 1).Value = "Caltrain Pilot"
$ws.Cells.Item(548, 2).Value = "TRANSFERS_FROM_CODE"
$ws.Cells.Item(548, 3).Value = 0
$ws.Cells.Item(548, 4).Value = "origin_survey_board_transfers"
$ws.Cells.Item(548, 5).Value = "none"

$ws.Cells.Item(549, 1).Value = "Caltrain Pilot"
$ws.Cells.Item(549, 2).Value = "TRANSFERS_FROM_CODE"
$ws.Cells.Item(549, 3).Value = 1
$ws.Cells.Item(549, 4).Value = "origin_survey_board_transfers"
$ws.Cells.Item(549, 5).Value = "one"

$ws.Cells.Item(550, 1).Value = "Caltrain Pilot"
$ws.Cells.Item(550, 2).Value = "TRANSFERS_FROM_CODE"
$ws.Cells.Item(550, 3).Value = 2
$ws.Cells.Item(550, 4).Value = "origin_survey_board_transfers"
$ws.Cells.Item(550, 5).Value = "two"

$ws.Cells.Item(551, 1).Value = "Caltrain Pilot"
$ws.Cells.Item(551, 2).Value = "TRANSFERS_FROM_CODE"
$ws.Cells.Item(551, 3).Value = "'3+"
$ws.Cells.Item(551, 4).Value = "origin_survey_board_transfers"
$ws.Cells.Item(551, 5).Value = "three or more"

$ws.Cells.Item(552, 1).Value = "Caltrain Pilot"
$ws.Cells.Item(552, 2).Value = "TRANSFERS_TO_CODE"
$ws.Cells.Item(552, 3).Value = 0
$ws.Cells.Item(552, 4).Value = "survey_alight_dest_transfers"
$ws.Cells.Item(552, 5).Value = "none"

$ws.Cells.Item(553, 1).Value = "Caltrain Pilot"
$ws.Cells.Item(553, 2).Value = "TRANSFERS_TO_CODE"
$ws.Cells.Item(553, 3).Value = 1
$ws.Cells.Item(553, 4).Value = "survey_alight_dest_transfers"
$ws.Cells.Item(553, 5).Value = "one"

$ws.Cells.Item(554, 1).Value = "Caltrain Pilot"
$ws.Cells.Item(554, 2).Value = "TRANSFERS_TO_CODE"
$ws.Cells.Item(554, 3).Value = 2
$ws.Cells.Item(554, 4).Value = "survey_alight_dest_transfers"
$ws.Cells.Item(554, 5).Value = "two"

$ws.Cells.Item(555, 1).Value = "Caltrain Pilot"
$ws.Cells.Item(555, 2).Value = "TRANSFERS_TO_CODE"
$ws.Cells.Item(555, 3).Value = "'3+"
$ws.Cells.Item(555, 4).Value = "survey_alight_dest_transfers"
$ws.Cells.Item(555, 5).Value = "three or more"

$ws.Cells.Item(556, 1).Value = "Caltrain Pilot"
$ws.Cells.Item(556, 2).Value = "TRANSFER_FROM_1ST"
$ws.Cells.Item(556, 3).Value = "NONCATEGORICAL"
$ws.Cells.Item(556, 4).Value = "first_route_before_survey_board"
$ws.Cells.Item(556, 5).Value = "NONCATEGORICAL"

$ws.Cells.Item(557, 1).Value = "Caltrain Pilot"
$ws.Cells.Item(557, 2).Value = "TRANSFER_FROM_2ND"
$ws.Cells.Item(557, 3).Value = "NONCATEGORICAL"
$ws.Cells.Item(557, 4).Value = "second_route_before_survey_board"
$ws.Cells.Item(557, 5).Value = "NONCATEGORICAL"

$ws.Cells.Item(558, 1).Value = "Caltrain Pilot"
$ws.Cells.Item(558, 2).Value = "TRANSFER_FROM_3RD"
$ws.Cells.Item(558, 3).Value = "NONCATEGORICAL"
$ws.Cells.Item(558, 4).Value = "third_route_before_survey_board"
$ws.Cells.Item(558, 5).Value = "NONCATEGORICAL"

$ws.Cells.Item(559, 1).Value = "Caltrain Pilot"
$ws.Cells.Item(559, 2).Value = "TRANSFER_TO_1ST"
$ws.Cells.Item(559, 3).Value = "NONCATEGORICAL"
$ws.Cells.Item(559, 4).Value = "first_route_after_survey_alight"
$ws.Cells.Item(559, 5).Value = "NONCATEGORICAL"

$ws.Cells.Item(560, 1).Value = "Caltrain Pilot"
$ws.Cells.Item(560, 2).Value = "TRANSFER_TO_2ND"
$ws.Cells.Item(560, 3).Value = "NONCATEGORICAL"
$ws.Cells.Item(560, 4).Value = "second_route_after_survey_alight"
$ws.Cells.Item(560, 5).Value = "NONCATEGORICAL"

$ws.Cells.Item(561, 1).Value = "Caltrain Pilot"
$ws.Cells.Item(561, 2).Value = "TRANSFER_TO_3RD"
$ws.Cells.Item(561, 3).Value = "NONCATEGORICAL"
$ws.Cells.Item(561, 4).Value = "third_route_after_survey_alight"
$ws.Cells.Item(561, 5).Value = "NONCATEGORICAL"

# Column D needs to widen to fit the new longer "origin_survey_board_transfers" /
# "survey_alight_dest_transfers" text (was 25.375, target 30.625)
$ws.Columns.Item(4).ColumnWidth = 29.83

# Move the frozen-pane view / selection down to the new bottom of the sheet
$ws.Range("A562").Select()
$excel.ActiveWindow.ScrollRow = 532
